# Auto-generated edit script applying odds updates from the diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").Value = 1.06
$ws.Range("K2").Value = 10
$ws.Range("L2").Value = 1.33
$ws.Range("M2").Value = 3.4
$ws.Range("N2").Value = 2
$ws.Range("O2").Value = 1.8

# Row 3
$ws.Range("J3").Value = 1.07
$ws.Range("K3").Value = 8.5
$ws.Range("R3").Value = 1.91
$ws.Range("S3").Value = 1.91
$ws.Range("Z3").Value = 8.5
$ws.Range("AD3").Value = 401

# Row 4
$ws.Range("G4").Value = 2.25
$ws.Range("I4").Value = 3.2
$ws.Range("N4").Value = 1.92
$ws.Range("O4").Value = 1.98
$ws.Range("R4").Value = 1.7
$ws.Range("S4").Value = 2.05
$ws.Range("T4").Value = 9.5
$ws.Range("U4").Value = 13
$ws.Range("X4").Value = 19
$ws.Range("Z4").Value = 11
$ws.Range("AD4").Value = 251
$ws.Range("AE4").Value = 12
$ws.Range("AF4").Value = 19
$ws.Range("AG4").Value = 13
$ws.Range("AH4").Value = 41

# Row 5
$ws.Range("G5").Value = 1.29
$ws.Range("H5").Value = 6
$ws.Range("I5").Value = 8.5
$ws.Range("N5").Value = 1.4
$ws.Range("O5").Value = 2.88
$ws.Range("P5").Value = 1.22
$ws.Range("Q5").Value = 4
$ws.Range("T5").Value = 10
$ws.Range("U5").Value = 8
$ws.Range("V5").Value = 9
$ws.Range("Z5").Value = 21
$ws.Range("AA5").Value = 12
$ws.Range("AB5").Value = 21

# Row 8
$ws.Range("G8").Value = 2.2
$ws.Range("I8").Value = 3.75
$ws.Range("U8").Value = 9
$ws.Range("V8").Value = 10
$ws.Range("W8").Value = 21
$ws.Range("AE8").Value = 8
$ws.Range("AF8").Value = 17
$ws.Range("AG8").Value = 15
$ws.Range("AI8").Value = 41

# Row 9
$ws.Range("G9").Value = 4.2
$ws.Range("H9").Value = 2.75
$ws.Range("J9").Value = 1.18
$ws.Range("K9").Value = 4.5
$ws.Range("L9").Value = 1.8
$ws.Range("M9").Value = 1.91
$ws.Range("N9").Value = 3.6
$ws.Range("O9").Value = 1.29
$ws.Range("P9").Value = 1.83
$ws.Range("Q9").Value = 1.98
$ws.Range("R9").Value = 2.75
$ws.Range("S9").Value = 1.4
$ws.Range("T9").Value = 7
$ws.Range("W9").Value = 51
$ws.Range("X9").Value = 51
$ws.Range("Z9").Value = 4.33
$ws.Range("AB9").Value = 26
$ws.Range("AC9").Value = 126
$ws.Range("AE9").Value = 4.75
$ws.Range("AG9").Value = 12
$ws.Range("AI9").Value = 29

# Row 10
$ws.Range("G10").Value = 1.8
$ws.Range("H10").Value = 3.6
$ws.Range("I10").Value = 4.1
$ws.Range("P10").Value = 1.33
$ws.Range("Q10").Value = 3.25
$ws.Range("X10").Value = 13
$ws.Range("AF10").Value = 23
$ws.Range("AI10").Value = 34

# Row 11
$ws.Range("G11").Value = 1.7
$ws.Range("H11").Value = 4
$ws.Range("J11").Value = 1.03
$ws.Range("K11").Value = 15
$ws.Range("L11").Value = 1.18
$ws.Range("M11").Value = 4.5
$ws.Range("N11").Value = 1.62
$ws.Range("O11").Value = 2.25
$ws.Range("P11").Value = 1.3
$ws.Range("Q11").Value = 3.4
$ws.Range("R11").Value = 1.62
$ws.Range("S11").Value = 2.2
$ws.Range("T11").Value = 9
$ws.Range("U11").Value = 9.5
$ws.Range("Y11").Value = 21
$ws.Range("Z11").Value = 15
$ws.Range("AA11").Value = 8
$ws.Range("AB11").Value = 13
$ws.Range("AH11").Value = 51

# Row 12
$ws.Range("G12").Value = 2.7
$ws.Range("H12").Value = 3.25
$ws.Range("I12").Value = 2.55
$ws.Range("J12").Value = 1.07
$ws.Range("K12").Value = 9
$ws.Range("V12").Value = 11
$ws.Range("W12").Value = 29

# Row 13
$ws.Range("K13").Value = 8

# Row 14
$ws.Range("G14").Value = 2.7
$ws.Range("I14").Value = 2.38
$ws.Range("K14").Value = 17
$ws.Range("N14").Value = 1.57
$ws.Range("O14").Value = 2.35
$ws.Range("R14").Value = 1.5
$ws.Range("S14").Value = 2.5
$ws.Range("X14").Value = 19
$ws.Range("Z14").Value = 17
$ws.Range("AB14").Value = 11
$ws.Range("AE14").Value = 12
$ws.Range("AF14").Value = 15

# Row 16
$ws.Range("G16").Value = 3.3
$ws.Range("I16").Value = 2.15
$ws.Range("R16").Value = 1.95
$ws.Range("S16").Value = 1.8
$ws.Range("U16").Value = 17
$ws.Range("Z16").Value = 9
$ws.Range("AD16").Value = 351
$ws.Range("AG16").Value = 9
$ws.Range("AH16").Value = 19

# Row 17
$ws.Range("G17").Value = 2.7
$ws.Range("I17").Value = 2.5
$ws.Range("J17").Value = 1.07
$ws.Range("K17").Value = 9
$ws.Range("V17").Value = 11
$ws.Range("W17").Value = 29
$ws.Range("X17").Value = 23
$ws.Range("Z17").Value = 9
$ws.Range("AD17").Value = 301
$ws.Range("AE17").Value = 8
$ws.Range("AF17").Value = 12
$ws.Range("AH17").Value = 23

# Row 18
$ws.Range("G18").Value = 1.83
$ws.Range("I18").Value = 4.33
$ws.Range("N18").Value = 2.4
$ws.Range("O18").Value = 1.53
$ws.Range("R18").Value = 2.2
$ws.Range("S18").Value = 1.62
$ws.Range("U18").Value = 7.5
$ws.Range("W18").Value = 15
$ws.Range("Y18").Value = 41
$ws.Range("Z18").Value = 7
$ws.Range("AB18").Value = 21
$ws.Range("AC18").Value = 81
$ws.Range("AF18").Value = 21
$ws.Range("AH18").Value = 51

# Row 19
$ws.Range("G19").Value = 3.4
$ws.Range("K19").Value = 8
$ws.Range("N19").Value = 2.35
$ws.Range("O19").Value = 1.57
$ws.Range("AD19").Value = 401
$ws.Range("AF19").Value = 9.5

# Row 20
$ws.Range("I20").Value = 9.5
$ws.Range("J20").Value = 1.07
$ws.Range("K20").Value = 8.5
$ws.Range("N20").Value = 2.08
$ws.Range("O20").Value = 1.73
$ws.Range("P20").Value = 1.44
$ws.Range("Q20").Value = 2.63
$ws.Range("R20").Value = 2.63
$ws.Range("S20").Value = 1.44
$ws.Range("T20").Value = 5
$ws.Range("W20").Value = 7.5
$ws.Range("X20").Value = 15
$ws.Range("Z20").Value = 8.5
$ws.Range("AC20").Value = 126
$ws.Range("AH20").Value = 151

# Row 21
$ws.Range("G21").Value = 7.5
$ws.Range("I21").Value = 1.38
$ws.Range("N21").Value = 1.85
$ws.Range("O21").Value = 2
$ws.Range("R21").Value = 2.1
$ws.Range("S21").Value = 1.67
$ws.Range("AA21").Value = 9
$ws.Range("AD21").Value = 501

# Row 22
$ws.Range("H22").Value = 3.6
$ws.Range("I22").Value = 2.2
$ws.Range("N22").Value = 1.95
$ws.Range("P22").Value = 1.4
$ws.Range("Q22").Value = 2.75
$ws.Range("R22").Value = 1.8
$ws.Range("S22").Value = 1.91
$ws.Range("T22").Value = 9.5
$ws.Range("Z22").Value = 11
$ws.Range("AA22").Value = 7
$ws.Range("AE22").Value = 8
$ws.Range("AG22").Value = 9
$ws.Range("AI22").Value = 17

# Row 23
$ws.Range("G23").Value = 2.9
$ws.Range("H23").Value = 3.4
$ws.Range("I23").Value = 2.25
$ws.Range("J23").Value = 1.05
$ws.Range("K23").Value = 11
$ws.Range("L23").Value = 1.3
$ws.Range("M23").Value = 3.4
$ws.Range("N23").Value = 2
$ws.Range("O23").Value = 1.85
$ws.Range("P23").Value = 1.4
$ws.Range("Q23").Value = 2.75
$ws.Range("R23").Value = 1.73
$ws.Range("S23").Value = 2
$ws.Range("T23").Value = 9.5
$ws.Range("W23").Value = 34
$ws.Range("X23").Value = 23
$ws.Range("Y23").Value = 34
$ws.Range("Z23").Value = 11
$ws.Range("AA23").Value = 6.5
$ws.Range("AB23").Value = 15
$ws.Range("AC23").Value = 51
$ws.Range("AE23").Value = 8
$ws.Range("AF23").Value = 11
$ws.Range("AJ23").Value = 29

# Row 24
$ws.Range("G24").Value = 3.9
$ws.Range("H24").Value = 3.25
$ws.Range("I24").Value = 2.05
$ws.Range("X24").Value = 41
$ws.Range("AA24").Value = 6.5

# Row 26
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 3.4
$ws.Range("I26").Value = 2
$ws.Range("R26").Value = 1.75
$ws.Range("S26").Value = 2
$ws.Range("V26").Value = 11
$ws.Range("W26").Value = 34

# Row 27
$ws.Range("N27").Value = 1.67
$ws.Range("O27").Value = 2.15

# Row 29
$ws.Range("H29").Value = 3.85
$ws.Range("I29").Value = 4.4
$ws.Range("L29").Value = 1.21
$ws.Range("M29").Value = 3.95
$ws.Range("N29").Value = 1.65
$ws.Range("O29").Value = 2.12
$ws.Range("P29").Value = 1.31
$ws.Range("Q29").Value = 3.15
$ws.Range("R29").Value = 1.65
$ws.Range("S29").Value = 2.12
$ws.Range("T29").Value = 8.5
$ws.Range("U29").Value = 9
$ws.Range("W29").Value = 13.5
$ws.Range("AA29").Value = 7.6
$ws.Range("AB29").Value = 14
$ws.Range("AC29").Value = 55
$ws.Range("AD29").Value = 350
$ws.Range("AE29").Value = 14.5
$ws.Range("AJ29").Value = 37
